# Update capital structure database values for Romania / Investments & Asset Management
# rows 2-7, per the source-data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = -0.3855
$ws.Range("E2").Value = -0.268
$ws.Range("G2").Value = 0.6248272048659109
$ws.Range("H2").Value = 0.6248272048659109
$ws.Range("I2").Value = 0.5340653264346933
$ws.Range("J2").Value = 0.4620313541173277
$ws.Range("K2").Value = 131.09
$ws.Range("L2").Value = 0.5177534657766895
$ws.Range("M2").Value = 355.96
$ws.Range("N2").Value = 0.1086237412267318
$ws.Range("O2").Value = 2.71538637577237
$ws.Range("P2").Value = 112.83
$ws.Range("Q2").Value = 0.03443088190418065
$ws.Range("R2").Value = 0.8607063849263864
$ws.Range("S2").Value = 243.13
$ws.Range("T2").Value = 0.6830261827171592
$ws.Range("U2").Value = 306.32
$ws.Range("V2").Value = 0.09347574000610315
$ws.Range("W2").Value = 0.02193744569939183
$ws.Range("X2").Value = 0.04699080201778191
$ws.Range("Y2").Value = -0.02505335631839008
$ws.Range("Z2").Value = 0.06038310923287527
$ws.Range("AA2").Value = 0.02065097581534072
$ws.Range("AB2").Value = 0.04697868048374997
$ws.Range("AC2").Value = -0.02693248796158438
$ws.Range("AD2").Value = 12.258
$ws.Range("AF2").Value = 12.258
$ws.Range("AG2").Value = -294.062
$ws.Range("AH2").Value = 0.003726676350714964
$ws.Range("AI2").Value = 0.003006727239457445
$ws.Range("AJ2").Value = -0.09858133155969048
$ws.Range("AK2").Value = -0.07798940098203493
$ws.Range("AL2").Value = 0.478
$ws.Range("AM2").Value = 0.478
$ws.Range("AN2").Value = 0.5281344248168892
$ws.Range("AO2").Value = 282.8870292887029
$ws.Range("AP2").Value = -12.66962516156829
$ws.Range("AQ2").Value = 282.8870292887029
# Row 3
$ws.Range("I3").Value = 0.846774193548387
$ws.Range("J3").Value = 0.8280877778089178
$ws.Range("K3").Value = 112.6
$ws.Range("L3").Value = 0.8255131964809383
$ws.Range("M3").Value = 301.3
$ws.Range("N3").Value = 0.1330947963601025
$ws.Range("O3").Value = 2.675843694493783
$ws.Range("P3").Value = 95.5
$ws.Range("Q3").Value = 0.04218570545101157
$ws.Range("R3").Value = 0.8481349911190054
$ws.Range("S3").Value = 205.8
$ws.Range("T3").Value = 0.6830401593096581
$ws.Range("U3").Value = 230.4
$ws.Range("V3").Value = 0.101775775245163
$ws.Range("W3").Value = 0.0456369310582418
$ws.Range("X3").Value = 0.04696437217823866
$ws.Range("Y3").Value = -0.001327441119996865
$ws.Range("Z3").Value = 0.05655526992287917
$ws.Range("AA3").Value = 0.04683272779382055
$ws.Range("AB3").Value = 0.04696437217823866
$ws.Range("AC3").Value = -0.0001316443844181156
$ws.Range("AG3").Value = -230.4
$ws.Range("AJ3").Value = -0.1133077604012983
$ws.Range("AK3").Value = -0.1068001668752608
$ws.Range("AL3").Value = 0.017
$ws.Range("AM3").Value = 0.017
$ws.Range("AO3").Value = 6794.117647058823
$ws.Range("AQ3").Value = 6794.117647058823
# Row 4
$ws.Range("D4").Value = -0.277
$ws.Range("E4").Value = -0.3
$ws.Range("G4").Value = 3.820754716981132
$ws.Range("H4").Value = 3.820754716981132
$ws.Range("I4").Value = 0.5688679245283019
$ws.Range("J4").Value = 0.513547929325232
$ws.Range("K4").Value = 8.52
$ws.Range("L4").Value = 0.8037735849056603
$ws.Range("M4").Value = 4.22
$ws.Range("N4").Value = 0.02257891920813269
$ws.Range("O4").Value = 0.4953051643192488
$ws.Range("P4").Value = 4.22
$ws.Range("Q4").Value = 0.02257891920813269
$ws.Range("R4").Value = 0.4953051643192488
$ws.Range("U4").Value = 8.35
$ws.Range("V4").Value = 0.04467629748528625
$ws.Range("W4").Value = 0.03243243243243243
$ws.Range("X4").Value = 0.04696437217823866
$ws.Range("Y4").Value = -0.01453193974580623
$ws.Range("Z4").Value = 0.04574881312041433
$ws.Range("AA4").Value = 0.02349420824707578
$ws.Range("AB4").Value = 0.04696437217823866
$ws.Range("AC4").Value = -0.02347016393116289
$ws.Range("AG4").Value = -8.35
$ws.Range("AJ4").Value = -0.0467656118734248
$ws.Range("AK4").Value = -0.03629645729189306
$ws.Range("AP4").Value = -1.353322528363047
# Row 5
$ws.Range("G5").Value = 0.589527027027027
$ws.Range("H5").Value = 0.589527027027027
$ws.Range("I5").Value = 0.1773648648648649
$ws.Range("J5").Value = 0.161402027027027
$ws.Range("K5").Value = 10.1
$ws.Range("L5").Value = 0.1706081081081081
$ws.Range("M5").Value = 31.41
$ws.Range("N5").Value = 0.1329805249788315
$ws.Range("O5").Value = 3.10990099009901
$ws.Range("P5").Value = 2.41
$ws.Range("Q5").Value = 0.01020321761219306
$ws.Range("R5").Value = 0.2386138613861386
$ws.Range("S5").Value = 29
$ws.Range("T5").Value = 0.9232728430436167
$ws.Range("U5").Value = 10.6
$ws.Range("V5").Value = 0.04487722269263336
$ws.Range("W5").Value = 0.02193744569939183
$ws.Range("X5").Value = 0.04814979530598221
$ws.Range("Y5").Value = -0.02621234960659037
$ws.Range("Z5").Value = 0.1279474378093324
$ws.Range("AA5").Value = 0.02065097581534072
$ws.Range("AB5").Value = 0.04758346377692511
$ws.Range("AC5").Value = -0.02693248796158438
$ws.Range("AD5").Value = 8.85
$ws.Range("AF5").Value = 8.85
$ws.Range("AG5").Value = -1.75
$ws.Range("AH5").Value = 0.03611507855539686
$ws.Range("AI5").Value = 0.02138972809667674
$ws.Range("AJ5").Value = -0.007464278097675411
$ws.Range("AK5").Value = -0.004340816073421803
$ws.Range("AN5").Value = 0.6807692307692308
$ws.Range("AP5").Value = -0.1346153846153846
# Row 6
$ws.Range("B6").Value = "Societatea de Investitii Financiare Moldova S.A. (BVB:SIF2)"
$ws.Range("E6").Value = -0.236
$ws.Range("G6").Value = 0.9452054794520548
$ws.Range("H6").Value = 0.9452054794520548
$ws.Range("I6").Value = 0.1321917808219178
$ws.Range("J6").Value = 0.0707121113285497
$ws.Range("K6").Value = 2.73
$ws.Range("L6").Value = 0.06232876712328767
$ws.Range("M6").Value = 17.73
$ws.Range("N6").Value = 0.05809305373525557
$ws.Range("O6").Value = 6.494505494505495
$ws.Range("P6").Value = 10.7
$ws.Range("Q6").Value = 0.03505897771952818
$ws.Range("R6").Value = 3.919413919413919
$ws.Range("S6").Value = 7.030000000000001
$ws.Range("T6").Value = 0.3965031020868585
$ws.Range("U6").Value = 2.77
$ws.Range("V6").Value = 0.009076015727391874
$ws.Range("W6").Value = 0.005681581685744017
$ws.Range("X6").Value = 0.04729298518687632
$ws.Range("Y6").Value = -0.0416114035011323
$ws.Range("Z6").Value = 0.09024787258154245
$ws.Range("AA6").Value = 0.006381617613150798
$ws.Range("AB6").Value = 0.04714059148487108
$ws.Range("AC6").Value = -0.04075897387172028
$ws.Range("AD6").Value = 3.17
$ws.Range("AF6").Value = 3.17
$ws.Range("AG6").Value = 0.3999999999999999
$ws.Range("AH6").Value = 0.01027985861140837
$ws.Range("AI6").Value = 0.00722145021299861
$ws.Range("AJ6").Value = 0.001308900523560209
$ws.Range("AK6").Value = 0.0009170105456212745
$ws.Range("AL6").Value = 0.435
$ws.Range("AM6").Value = 0.435
$ws.Range("AN6").Value = 0.483969465648855
$ws.Range("AO6").Value = 13.31034482758621
$ws.Range("AP6").Value = 0.06106870229007633
$ws.Range("AQ6").Value = 13.31034482758621
# Row 7
$ws.Range("B7").Value = "SIF Banat-Crisana (BVB:SIF1)"
$ws.Range("D7").Value = -0.494
$ws.Range("G7").Value = 12.97805642633229
$ws.Range("H7").Value = 12.97805642633229
$ws.Range("I7").Value = -0.8150470219435737
$ws.Range("J7").Value = -0.8150470219435737
$ws.Range("K7").Value = -2.86
$ws.Range("L7").Value = -0.896551724137931
$ws.Range("M7").Value = 1.3
$ws.Range("N7").Value = 0.004563004563004564
$ws.Range("O7").Value = -0.4545454545454546
$ws.Range("P7").Value = -0
$ws.Range("Q7").Value = -0
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 1.3
$ws.Range("T7").Value = 1
$ws.Range("U7").Value = 54.2
$ws.Range("V7").Value = 0.1902421902421903
$ws.Range("W7").Value = -0.004570880613712642
$ws.Range("X7").Value = 0.04699080201778191
$ws.Range("Y7").Value = -0.05156168263149456
$ws.Range("Z7").Value = 0.005303055490906672
$ws.Range("AA7").Value = -0.004322239585064999
$ws.Range("AB7").Value = 0.04697868048374997
$ws.Range("AC7").Value = -0.05130092006881497
$ws.Range("AD7").Value = 0.238
$ws.Range("AF7").Value = 0.238
$ws.Range("AG7").Value = -53.962
$ws.Range("AH7").Value = 0.0008346835567339324
$ws.Range("AI7").Value = 0.0003979680220989301
$ws.Range("AJ7").Value = -0.2336644467346214
$ws.Range("AK7").Value = -0.09922440138423576
$ws.Range("AL7").Value = 0.026
$ws.Range("AM7").Value = 0.026
$ws.Range("AN7").Value = -0.09482071713147411
$ws.Range("AO7").Value = -100
$ws.Range("AP7").Value = 21.4988047808765
$ws.Range("AQ7").Value = -100

# Cells removed entirely by this edit (now blank)
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("E7").ClearContents()
